# Weekly data refresh: insert a new daily price record as the new row 740,
# pushing all subsequent rows down by one (sheet grows from 790 to 791 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 740; existing rows 740.. shift down to 741..
$ws.Rows(740).Insert()

# The newly inserted row 740 is blank. Populate it. Most of the "dimension"
# columns (market/region/product taxonomy) are identical for every row in
# this sheet, so copy them from the row directly below (the old row 740,
# now shifted to 741) and only set the columns that carry the new record's
# actual data.
$ws.Cells.Item(740, 1).Value2  = $ws.Cells.Item(741, 1).Value2   # Mercado ID
$ws.Cells.Item(740, 2).Value2  = $ws.Cells.Item(741, 2).Value2   # Mercado
$ws.Cells.Item(740, 3).Value2  = $ws.Cells.Item(741, 3).Value2   # Región
$ws.Cells.Item(740, 4).Value2  = 44939                           # Fecha (2023-01-13)
$ws.Cells.Item(740, 5).Value2  = $ws.Cells.Item(741, 5).Value2   # Codreg
$ws.Cells.Item(740, 6).Value2  = $ws.Cells.Item(741, 6).Value2   # Tipo
$ws.Cells.Item(740, 7).Value2  = $ws.Cells.Item(741, 7).Value2   # Producto ID
$ws.Cells.Item(740, 8).Value2  = $ws.Cells.Item(741, 8).Value2   # Producto
$ws.Cells.Item(740, 9).Value2  = $ws.Cells.Item(741, 9).Value2   # Categoría ID
$ws.Cells.Item(740, 10).Value2 = $ws.Cells.Item(741, 10).Value2  # Categoría
$ws.Cells.Item(740, 11).Value2 = "Murcott"                       # Variedad
$ws.Cells.Item(740, 12).Value2 = "Primera"                       # Calidad
$ws.Cells.Item(740, 13).Value2 = 80                              # Volumen
$ws.Cells.Item(740, 14).Value2 = 12000                           # Precio mínimo
$ws.Cells.Item(740, 15).Value2 = 12000                           # Precio máximo
$ws.Cells.Item(740, 16).Value2 = 12000                           # Precio promedio ponderado
$ws.Cells.Item(740, 17).Value2 = "$/bandeja 10 kilos"            # Unidad de comercialización
$ws.Cells.Item(740, 18).Value2 = "Región de O'Higgins"           # Origen
$ws.Cells.Item(740, 19).Value2 = 1200                            # Precio $/Kg
$ws.Cells.Item(740, 20).Value2 = 10                              # Kg / unidad

# Match the date-cell number format used by the rest of column D.
$ws.Cells.Item(740, 4).NumberFormat = $ws.Cells.Item(741, 4).NumberFormat
